# Commit: "Update posts.xlsx after post"
# The post formerly on row 583 ("花より団子おおお（違う）...") was removed.
# Deleting its entire row shifts every subsequent row up by one
# (old row 584 -> new row 583, ..., old row 729 -> new row 728),
# and the sheet's used range shrinks from A1:C729 to A1:C728.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(583).Delete()
